$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "1" to "Khelvachauri"
$ws.Name = "Khelvachauri"

# Row 6 ("Urban") - mark every year (2010-2023, columns B:O) as unavailable/confidential
for ($col = 2; $col -le 15; $col++) {
    $ws.Cells.Item(6, $col).Value2 = "..."
}

# Row 7 ("Rural") - only the 2010 figure (column B) becomes unavailable/confidential
$ws.Cells.Item(7, 2).Value2 = "..."

# Remove the now-empty row 8, shifting the footnote row (old row 9) up to row 8
$ws.Rows(8).Delete()
